$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2330.476
$ws.Range("J17").Value = 2330.476
$ws.Range("L17").Value = 6991.428
$ws.Range("N17").Value = -7327.428
$ws.Range("H20").Value = 696.5
$ws.Range("I20").Value = 696.5
$ws.Range("K20").Value = 696.5
$ws.Range("M20").Value = -466.5
$ws.Range("H35").Value = 696.5
$ws.Range("I35").Value = 696.5
$ws.Range("K35").Value = 696.5
$ws.Range("M35").Value = -317.5
$ws.Range("H43").Value = 850
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -1138
$ws.Range("H100").Value = 718.93335
$ws.Range("I100").Value = 750.7143
$ws.Range("K100").Value = 750.7143
$ws.Range("M100").Value = -209.7143
$ws.Range("H106").Value = 4166
$ws.Range("I106").Value = 4166
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4166
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3535
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 1502.4193
$ws.Range("I132").Value = 1502.4193
$ws.Range("K132").Value = 4507.257900000001
$ws.Range("M132").Value = -1977.257900000001
$ws.Range("H138").Value = 2928.5715
$ws.Range("J138").Value = 3333.3333
$ws.Range("L138").Value = 9999.999899999999
$ws.Range("N138").Value = -20279.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4067.0967
$ws.Range("I32").Value = 3881.5862
$ws.Range("K32").Value = 3881.5862
$ws.Range("M32").Value = -3594.5862
$ws.Range("H45").Value = 1980.4445
$ws.Range("I45").Value = 1915.5
$ws.Range("K45").Value = 1915.5
$ws.Range("M45").Value = -1538.5
$ws.Range("H46").Value = 5399.6665
$ws.Range("J46").Value = 6200
$ws.Range("L46").Value = 6200
$ws.Range("N46").Value = -6838
$ws.Range("H61").Value = 2605.4167
$ws.Range("I61").Value = 2660.4546
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2660.4546
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2448.4546
$ws.Range("N61").Value = -2424
$ws.Range("H97").Value = 598.6
$ws.Range("I97").Value = 598.6
$ws.Range("K97").Value = 598.6
$ws.Range("M97").Value = -102.6
$ws.Range("H106").Value = 49000
$ws.Range("J106").Value = 49000
$ws.Range("L106").Value = 49000
$ws.Range("N106").Value = -51524
$ws.Range("H132").Value = 1581.7778
$ws.Range("I132").Value = 1557.1765
$ws.Range("K132").Value = 4671.529500000001
$ws.Range("M132").Value = -2141.529500000001
$ws.Range("H136").Value = 2605.4167
$ws.Range("I136").Value = 2660.4546
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 7981.3638
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -5431.3638
$ws.Range("N136").Value = -11100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3236.2856
$ws.Range("I5").Value = 1335
$ws.Range("K5").Value = 1335
$ws.Range("M5").Value = -1222
$ws.Range("H7").Value = 275
$ws.Range("J7").Value = 312.5
$ws.Range("L7").Value = 312.5
$ws.Range("N7").Value = -538.5
$ws.Range("H105").Value = 2562.875
$ws.Range("I105").Value = 2000.6
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 2000.6
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -253.5999999999999
$ws.Range("N105").Value = -6994
$ws.Range("H134").Value = 3374.1538
$ws.Range("I134").Value = 3374.1538
$ws.Range("K134").Value = 10122.4614
$ws.Range("M134").Value = -7587.4614

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 774.75
$ws.Range("I2").Value = 549.5
$ws.Range("K2").Value = 549.5
$ws.Range("M2").Value = -436.5
$ws.Range("H7").Value = 45.083332
$ws.Range("I7").Value = 48
$ws.Range("K7").Value = 48
$ws.Range("M7").Value = 65
$ws.Range("H14").Value = 934.8333
$ws.Range("I14").Value = 305
$ws.Range("J14").Value = 1564.6666
$ws.Range("K14").Value = 305
$ws.Range("L14").Value = 1564.6666
$ws.Range("M14").Value = -135
$ws.Range("N14").Value = -1904.6666
$ws.Range("H43").Value = 10899.714
$ws.Range("J43").Value = 10899.714
$ws.Range("L43").Value = 10899.714
$ws.Range("N43").Value = -11267.714
$ws.Range("H101").Value = 10899.714
$ws.Range("J101").Value = 10899.714
$ws.Range("L101").Value = 10899.714
$ws.Range("N101").Value = -17389.714
$ws.Range("H122").Value = 626.8182
$ws.Range("I122").Value = 609.5
$ws.Range("K122").Value = 1828.5
$ws.Range("M122").Value = 621.5
$ws.Range("H132").Value = 902.1579
$ws.Range("I132").Value = 929.2778
$ws.Range("J132").Value = 414
$ws.Range("K132").Value = 2787.8334
$ws.Range("L132").Value = 1242
$ws.Range("M132").Value = -257.8334
$ws.Range("N132").Value = -6302

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9983.333000000001
$ws.Range("I43").Value = 4966.6665
$ws.Range("K43").Value = 4966.6665
$ws.Range("M43").Value = -4815.6665
$ws.Range("H46").Value = 15008.2
$ws.Range("I46").Value = 7520.5
$ws.Range("K46").Value = 7520.5
$ws.Range("M46").Value = -7364.5
$ws.Range("H122").Value = 3096.5
$ws.Range("I122").Value = 2084.5
$ws.Range("K122").Value = 6253.5
$ws.Range("M122").Value = -3803.5
$ws.Range("H126").Value = 4281.25
$ws.Range("I126").Value = 4178.5713
$ws.Range("K126").Value = 12535.7139
$ws.Range("M126").Value = -10065.7139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999
$ws.Range("I7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("M7").Value = -1887
$ws.Range("H16").Value = 825
$ws.Range("I16").Value = 305
$ws.Range("K16").Value = 305
$ws.Range("M16").Value = -135
$ws.Range("H93").Value = 962
$ws.Range("I93").Value = 946
$ws.Range("J93").Value = 994
$ws.Range("K93").Value = 946
$ws.Range("L93").Value = 994
$ws.Range("M93").Value = 302
$ws.Range("N93").Value = -3490
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527
$ws.Range("H132").Value = 4858.3335
$ws.Range("I132").Value = 2466.6667
$ws.Range("K132").Value = 7400.000100000001
$ws.Range("M132").Value = -4870.000100000001
$ws.Range("H136").Value = 3146.7
$ws.Range("I136").Value = 3146.7
$ws.Range("K136").Value = 9440.099999999999
$ws.Range("M136").Value = -6890.099999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 36666.332
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 36666.332
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 36666.332
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -37446.332
$ws.Range("H62").Value = 4740.4
$ws.Range("I62").Value = 4966.3335
$ws.Range("J62").Value = 4401.5
$ws.Range("K62").Value = 4966.3335
$ws.Range("L62").Value = 4401.5
$ws.Range("M62").Value = -4342.3335
$ws.Range("N62").Value = -5649.5
$ws.Range("H65").Value = 4740.4
$ws.Range("I65").Value = 4966.3335
$ws.Range("J65").Value = 4401.5
$ws.Range("K65").Value = 24831.6675
$ws.Range("L65").Value = 22007.5
$ws.Range("M65").Value = -21711.6675
$ws.Range("N65").Value = -28247.5
$ws.Range("H74").Value = 17675.834
$ws.Range("I74").Value = 16624
$ws.Range("J74").Value = 17886.2
$ws.Range("K74").Value = 16624
$ws.Range("L74").Value = 17886.2
$ws.Range("M74").Value = -15688
$ws.Range("N74").Value = -19758.2
$ws.Range("H77").Value = 17675.834
$ws.Range("I77").Value = 16624
$ws.Range("J77").Value = 17886.2
$ws.Range("K77").Value = 49872
$ws.Range("L77").Value = 53658.60000000001
$ws.Range("M77").Value = -45192
$ws.Range("N77").Value = -63018.60000000001
$ws.Range("H113").Value = 7297.7334
$ws.Range("I113").Value = 11528.223
$ws.Range("J113").Value = 952
$ws.Range("K113").Value = 34584.669
$ws.Range("L113").Value = 2856
$ws.Range("M113").Value = -32414.669
$ws.Range("N113").Value = -7196
$ws.Range("H132").Value = 2413.5715
$ws.Range("I132").Value = 2413.5715
$ws.Range("K132").Value = 7240.7145
$ws.Range("M132").Value = -4710.7145
$ws.Range("H136").Value = 2225.923
$ws.Range("I136").Value = 1684.2
$ws.Range("K136").Value = 5052.6
$ws.Range("M136").Value = -2502.6
